$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.614.37"
$ws.Range("E2").Value = "  -1.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.092.58"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "342.79"
$ws.Range("E5").Value = "  -1.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.009"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5171"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4361"
$ws.Range("E8").Value = "  -1.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.80"
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09206"
$ws.Range("E10").Value = "  +2.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.162"
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.61"
$ws.Range("E12").Value = "  -2.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.097.69"
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.727"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.170"
$ws.Range("E15").Value = "  +2.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.71"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001146"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.010"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.74"
$ws.Range("E19").Value = "  +1.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06669"
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.008"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.173"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.651.61"
$ws.Range("E23").Value = "  -1.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.36"
$ws.Range("E24").Value = "  -2.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.310"
$ws.Range("E25").Value = "  -2.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.349.18"
$ws.Range("E26").Value = "  -0.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.89"
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.54"
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.492"
$ws.Range("E29").Value = "  -1.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.06"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.127"
$ws.Range("E31").Value = "  -3.68%  "
$ws.Range("E32").Value = "  -1.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.622"
$ws.Range("E33").Value = "  -2.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.159"
$ws.Range("E34").Value = "  -1.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.965"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.261"
$ws.Range("E36").Value = "  +5.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.22"
$ws.Range("E37").Value = "  +0.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02569"
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06687"
$ws.Range("E39").Value = "  -1.76%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6992"
$ws.Range("E40").Value = "  +2.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.46"
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2218"
$ws.Range("E42").Value = "  -3.82%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.311"
$ws.Range("E43").Value = "  +4.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6885"
$ws.Range("E44").Value = "  +7.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.19"
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.290"
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.611"
$ws.Range("E47").Value = "  -1.76%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000346"
$ws.Range("E48").Value = "  -4.12%  "
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.216"
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "81.25"
$ws.Range("E50").Value = "  -1.72%  "
$ws.Range("B51").Value = "WEMIXTOKEN"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.170"
$ws.Range("E51").Value = "  -0.80%  "
